$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (B,C,D,E) stay as text, matching source data (inline strings)
# Row 36/37 swap: ARBITRUM and LidoDAOToken exchange rows, with updated D/E values
$ws.Range("B36:E37").NumberFormat = "@"
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.94"
$ws.Range("E36").Value = "  +4.76%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "3.05"
$ws.Range("E37").Value = "  -7.86%  "

# Remaining per-cell value updates (force text format first to preserve original text typing)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.007.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.247.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.02"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.63"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.89"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0829"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.53"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.589.96"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.853"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.31"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.250.02"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.890.25"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.44"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0978"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.38"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.81"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.41%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.64"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.85"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.43%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.13"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.80"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.12"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0844"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +9.13%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.49"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +18.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.68"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.17"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0314"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.763.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "74.60"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "81.17"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.16"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.72"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.68"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "57.20"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.61%  "
